$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 1.170631
$ws.Range("N2").Value = 2.341262
$ws.Range("O2").Value = 0.1596624636318675
$ws.Range("P2").Value = 0.1281983861842902
$ws.Range("Q2").Value = 0.08279365789566667
$ws.Range("R2").Value = 0.496761947374
$ws.Range("S2").Value = 0.1596624636318675
$ws.Range("T2").Value = 0.1281983861842902

# Row 3
$ws.Range("O3").Value = 0.4101137613801331
$ws.Range("P3").Value = 0.4939412918191532
$ws.Range("S3").Value = 0.4101137613801331
$ws.Range("T3").Value = 0.4939412918191532

# Row 4
$ws.Range("M4").Value = 0.1597873333333333
$ws.Range("N4").Value = 0.479362
$ws.Range("O4").Value = 0.02179340825346879
$ws.Range("P4").Value = 0.02624799565280337
$ws.Range("Q4").Value = 0.01130106567488889
$ws.Range("R4").Value = 0.101709591074
$ws.Range("S4").Value = 0.02179340825346879
$ws.Range("T4").Value = 0.02624799565280337

# Row 5
$ws.Range("M5").Value = 2.5622985
$ws.Range("N5").Value = 5.124597
$ws.Range("O5").Value = 0.349472114671693
$ws.Range("P5").Value = 0.2806029676494365
$ws.Range("Q5").Value = 0.1812202696115
$ws.Range("R5").Value = 1.087321617669
$ws.Range("S5").Value = 0.349472114671693
$ws.Range("T5").Value = 0.2806029676494365

# Row 6
$ws.Range("M6").Value = 0.2687716666666666
$ws.Range("N6").Value = 0.8063149999999999
$ws.Range("O6").Value = 0.03665779093022745
$ws.Range("P6").Value = 0.04415066821064279
$ws.Range("Q6").Value = 0.01900905530611111
$ws.Range("R6").Value = 0.171081497755
$ws.Range("S6").Value = 0.03665779093022745
$ws.Range("T6").Value = 0.04415066821064279

# Row 7
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.163505
$ws.Range("N7").Value = 0.490515
$ws.Range("O7").Value = 0.02230046113261011
$ws.Range("P7").Value = 0.02685869048367381
$ws.Range("Q7").Value = 0.01156400012833333
$ws.Range("R7").Value = 0.104076001155
$ws.Range("S7").Value = 0.02230046113261011
$ws.Range("T7").Value = 0.02685869048367381
